$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The value that used to live in C43 ("NA") now belongs to the new row 44,
# so C43 becomes empty.
$ws.Range("C43").Value = ""

# Append the new row 44.
$ws.Range("A44").Value = "'2025-04-09"
$ws.Range("A44").Style = "Normal"
$ws.Range("B44").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C44").Value = "NA"
$ws.Range("D44").Value = 1
